$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are text (prices formatted with dots as thousands separators).
# Force text type so Excel does not auto-convert single-dot values (e.g. "1.008") to numbers,
# then restore the default "Normal" style so no stray number-format style is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.868.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.51%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.60%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.51%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4602'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3691'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07172'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8770'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07848'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.862.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.322'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.385'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008706'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.889.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.992'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.972'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.967'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.917'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08798'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.125'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7556'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.465'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.131'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.559'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.083'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01929'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.927'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05115'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.924'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4964'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1594'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.280'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4673'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.007'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.606'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06103'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.74%  '
